$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '24.884.21'
$ws.Range("E2").Value = '  -1.01%  '

$ws.Range("D3").Value = "'" + '1.700.50'
$ws.Range("E3").Value = '  -1.17%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = "'" + '316.17'
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").Value = "'" + '1.003'
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = "'" + '0.4026'
$ws.Range("E7").Value = '  +1.18%  '

$ws.Range("D8").Value = "'" + '0.4055'
$ws.Range("E8").Value = '  -1.78%  '

$ws.Range("B9").Value = 'BinanceUSD'
$ws.Range("C9").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D9").Value = "'" + '1.004'
$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = "'" + '1.473'
$ws.Range("E10").Value = '  -4.08%  '

$ws.Range("D11").Value = "'" + '53.68'
$ws.Range("E11").Value = '  +0.11%  '

$ws.Range("D12").Value = "'" + '0.08796'
$ws.Range("E12").Value = '  -2.10%  '

$ws.Range("D13").Value = "'" + '26.00'
$ws.Range("E13").Value = '  +3.65%  '

$ws.Range("D14").Value = "'" + '7.458'
$ws.Range("E14").Value = '  -4.10%  '

$ws.Range("D15").Value = "'" + '8.035'
$ws.Range("E15").Value = '  -2.35%  '

$ws.Range("D16").Value = "'" + '0.00001347'
$ws.Range("E16").Value = '  -3.28%  '

$ws.Range("D17").Value = "'" + '1.695.96'
$ws.Range("E17").Value = '  -0.97%  '

$ws.Range("D18").Value = "'" + '95.98'
$ws.Range("E18").Value = '  -5.21%  '

$ws.Range("D19").Value = "'" + '0.07141'
$ws.Range("E19").Value = '  -0.31%  '

$ws.Range("D20").Value = "'" + '20.91'
$ws.Range("E20").Value = '  +2.70%  '

$ws.Range("D21").Value = "'" + '7.223'
$ws.Range("E21").Value = '  -3.65%  '

$ws.Range("D22").Value = "'" + '1.003'
$ws.Range("E22").Value = '  -0.25%  '

$ws.Range("D23").Value = "'" + '14.40'
$ws.Range("E23").Value = '  -1.39%  '

$ws.Range("D24").Value = "'" + '24.886.95'
$ws.Range("E24").Value = '  -1.02%  '

$ws.Range("D25").Value = "'" + '2.331'
$ws.Range("E25").Value = '  -0.74%  '

$ws.Range("D26").Value = "'" + '2.882'
$ws.Range("E26").Value = '  -8.01%  '

$ws.Range("B27").Value = 'HuobiToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D27").Value = "'" + '6.413'
$ws.Range("E27").Value = '  +23.05%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'" + '23.05'
$ws.Range("E28").Value = '  -0.59%  '

$ws.Range("D29").Value = "'" + '165.42'
$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").Value = "'" + '144.63'
$ws.Range("E30").Value = '  +3.02%  '

$ws.Range("D31").Value = "'" + '8.216'
$ws.Range("E31").Value = '  -11.02%  '

$ws.Range("D32").Value = "'" + '1.905.54'
$ws.Range("E32").Value = '  +0.41%  '

$ws.Range("D33").Value = "'" + '2.239'
$ws.Range("E33").Value = '  +13.74%  '

$ws.Range("D34").Value = "'" + '0.08920'
$ws.Range("E34").Value = '  -1.52%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = "'" + '7.362'
$ws.Range("E35").Value = '  -7.06%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = "'" + '0.03195'
$ws.Range("E36").Value = '  +7.16%  '

$ws.Range("D37").Value = "'" + '1.012'
$ws.Range("E37").Value = '  -7.40%  '

$ws.Range("D38").Value = "'" + '0.2835'
$ws.Range("E38").Value = '  +0.49%  '

$ws.Range("D39").Value = "'" + '0.8399'
$ws.Range("E39").Value = '  +2.34%  '

$ws.Range("D40").Value = "'" + '10.80'
$ws.Range("E40").Value = '  -2.93%  '

$ws.Range("D41").Value = "'" + '0.09329'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").Value = "'" + '14.06'
$ws.Range("E42").Value = '  -4.87%  '

$ws.Range("D43").Value = "'" + '1.466'
$ws.Range("E43").Value = '  -1.73%  '

$ws.Range("D44").Value = "'" + '17.44'
$ws.Range("E44").Value = '  +3.11%  '

$ws.Range("D45").Value = "'" + '2.701'
$ws.Range("E45").Value = '  +1.23%  '

$ws.Range("D46").Value = "'" + '0.7429'
$ws.Range("E46").Value = '  -0.28%  '

$ws.Range("D47").Value = "'" + '4.244'
$ws.Range("E47").Value = '  -1.13%  '

$ws.Range("D48").Value = "'" + '1.384'
$ws.Range("E48").Value = '  +1.21%  '

$ws.Range("D49").Value = "'" + '1.003'
$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("D50").Value = "'" + '141.76'
$ws.Range("E50").Value = '  +0.36%  '

$ws.Range("D51").Value = "'" + '0.08337'
$ws.Range("E51").Value = '  +2.63%  '

